$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.114.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.02%  '

$ws.Range("D3").Value = "'3.586.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'241.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.09%  '

$ws.Range("D6").Value = "'655.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.89%  '

$ws.Range("E7").Value = '  +14.96%  '

$ws.Range("D8").Value = "'0.422"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.99%  '

$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("E10").Value = '  +4.44%  '

$ws.Range("D11").Value = "'3.579.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").Value = "'44.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.80%  '

$ws.Range("E13").Value = '  +0.69%  '

$ws.Range("D14").Value = "'6.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.59%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = "'4.250.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.66%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = "'97.013.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.01%  '

$ws.Range("D17").Value = "'0.0000259"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.71%  '

$ws.Range("D18").Value = "'8.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +11.29%  '

$ws.Range("D19").Value = "'3.579.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.38%  '

$ws.Range("D20").Value = "'12.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.52%  '

$ws.Range("D21").Value = "'17.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '

$ws.Range("E22").Value = '  +8.88%  '

$ws.Range("D23").Value = "'3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("D24").Value = "'513.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.50%  '

$ws.Range("E25").Value = '  +5.17%  '

$ws.Range("D26").Value = "'6.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.21%  '

$ws.Range("D27").Value = "'100.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.86%  '

$ws.Range("D28").Value = "'12.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.15%  '

$ws.Range("D29").Value = "'3.775.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.63%  '

$ws.Range("D30").Value = "'0.158"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.51%  '

$ws.Range("E31").Value = '  -0.97%  '

$ws.Range("D32").Value = "'11.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.10%  '

$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("E34").Value = '  +3.34%  '

$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("E36").Value = '  -0.28%  '

$ws.Range("D37").Value = "'622.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.12%  '

$ws.Range("E38").Value = '  +4.07%  '

$ws.Range("D39").Value = "'0.564"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.11%  '

$ws.Range("E40").Value = '  +2.74%  '

$ws.Range("D41").Value = "'1.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.47%  '

$ws.Range("D42").Value = "'0.154"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.71%  '

$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("D44").Value = "'0.923"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.43%  '

$ws.Range("D45").Value = "'5.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.04%  '

$ws.Range("E46").Value = '  +6.22%  '

$ws.Range("D47").Value = "'2.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.59%  '

$ws.Range("D48").Value = "'23.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("D49").Value = "'33.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.36%  '

$ws.Range("D50").Value = "'8.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.59%  '

$ws.Range("D51").Value = "'3.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.11%  '
